$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values must stay text -- force Text format before writing
# so numeric-looking strings (e.g. "0.9998") are not auto-coerced to numbers
# and formatting (trailing zeros, grouping dots) is preserved verbatim.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.432.09'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.866.34'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7075'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.39'
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07871'
$ws.Range("E8").Value = '  -1.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3129'
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.51'
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07996'
$ws.Range("E11").Value = '  -3.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.927.70'
$ws.Range("E12").Value = '  +2.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.203'
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.45'
$ws.Range("E14").Value = '  -1.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6993'
$ws.Range("E15").Value = '  -2.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.466'
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008381'
$ws.Range("E17").Value = '  -2.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.580.74'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.63'
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.139.11'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.09'
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.622'
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1562'
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.028'
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.90'
$ws.Range("E27").Value = '  -1.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.73'
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("E29").Value = '  -0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.325'
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.276'
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.213'
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05307'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.890'
$ws.Range("E34").Value = '  -2.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7519'
$ws.Range("E35").Value = '  -2.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.170'
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.707'
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01881'
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.280.08'
$ws.Range("E39").Value = '  +1.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.742'
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8966'
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.61'
$ws.Range("E44").Value = '  -3.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("E46").Value = '  -3.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.021.00'
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.583'
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.794'
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.5177'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4308'
$ws.Range("E51").Value = '  -1.44%  '

# Rows 42/43: Quant and FraxShare swap places (coin + link + price + volume)
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.072'
$ws.Range("E42").Value = '  -6.49%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '108.85'
$ws.Range("E43").Value = '  -3.85%  '
